$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 7487.8965
$ws.Range("J69").Value = 7487.8965
$ws.Range("L69").Value = 22463.6895
$ws.Range("N69").Value = -24211.6895
$ws.Range("H72").Value = 7487.8965
$ws.Range("J72").Value = 7487.8965
$ws.Range("L72").Value = 67391.06849999999
$ws.Range("N72").Value = -76127.06849999999
$ws.Range("H137").Value = 2068
$ws.Range("I137").Value = 1202.8462
$ws.Range("J137").Value = 3473.875
$ws.Range("K137").Value = 3608.5386
$ws.Range("L137").Value = 10421.625
$ws.Range("M137").Value = -1058.5386
$ws.Range("N137").Value = -15521.625
$ws.Range("H138").Value = 4469.727
$ws.Range("I138").Value = 2645.875
$ws.Range("K138").Value = 7937.625
$ws.Range("M138").Value = -2797.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7693095
$ws.Range("I32").Value = 853.1667
$ws.Range("K32").Value = 853.1667
$ws.Range("M32").Value = -566.1667
$ws.Range("H45").Value = 2849.1667
$ws.Range("I45").Value = 2149.5557
$ws.Range("K45").Value = 2149.5557
$ws.Range("M45").Value = -1772.5557
$ws.Range("H74").Value = 6547.5
$ws.Range("I74").Value = 5000
$ws.Range("J74").Value = 8095
$ws.Range("K74").Value = 5000
$ws.Range("L74").Value = 8095
$ws.Range("M74").Value = -4126
$ws.Range("N74").Value = -9843
$ws.Range("H77").Value = 6547.5
$ws.Range("I77").Value = 5000
$ws.Range("J77").Value = 8095
$ws.Range("K77").Value = 25000
$ws.Range("L77").Value = 40475
$ws.Range("M77").Value = -20632
$ws.Range("N77").Value = -49211

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7422.4287
$ws.Range("I86").Value = 3975
$ws.Range("J86").Value = 8801.4
$ws.Range("K86").Value = 3975
$ws.Range("L86").Value = 8801.4
$ws.Range("M86").Value = -2852
$ws.Range("N86").Value = -11047.4
$ws.Range("H89").Value = 7422.4287
$ws.Range("I89").Value = 3975
$ws.Range("J89").Value = 8801.4
$ws.Range("K89").Value = 19875
$ws.Range("L89").Value = 44007
$ws.Range("M89").Value = -14259
$ws.Range("N89").Value = -55239
$ws.Range("H134").Value = 5622.647
$ws.Range("I134").Value = 2286.4666
$ws.Range("K134").Value = 6859.399800000001
$ws.Range("M134").Value = -4324.399800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6099.706
$ws.Range("I31").Value = 3973.75
$ws.Range("J31").Value = 6753.846
$ws.Range("K31").Value = 3973.75
$ws.Range("L31").Value = 6753.846
$ws.Range("M31").Value = -3678.75
$ws.Range("N31").Value = -7343.846
$ws.Range("H33").Value = 1232.6923
$ws.Range("I33").Value = 772.6
$ws.Range("J33").Value = 2766.3333
$ws.Range("K33").Value = 772.6
$ws.Range("L33").Value = 2766.3333
$ws.Range("M33").Value = -393.6
$ws.Range("N33").Value = -3524.3333
$ws.Range("H34").Value = 6099.706
$ws.Range("I34").Value = 3973.75
$ws.Range("J34").Value = 6753.846
$ws.Range("K34").Value = 3973.75
$ws.Range("L34").Value = 6753.846
$ws.Range("M34").Value = -3771.75
$ws.Range("N34").Value = -7157.846
$ws.Range("H62").Value = 1100
$ws.Range("I62").Value = 1100
$ws.Range("K62").Value = 1100
$ws.Range("M62").Value = -476
$ws.Range("H65").Value = 1100
$ws.Range("I65").Value = 1100
$ws.Range("K65").Value = 5500
$ws.Range("M65").Value = -2380
$ws.Range("H132").Value = 2332.158
$ws.Range("I132").Value = 1847.75
$ws.Range("J132").Value = 4915.6665
$ws.Range("K132").Value = 5543.25
$ws.Range("L132").Value = 14746.9995
$ws.Range("M132").Value = -3013.25
$ws.Range("N132").Value = -19806.9995
$ws.Range("H134").Value = 5648.5625
$ws.Range("I134").Value = 5108.615
$ws.Range("K134").Value = 15325.845
$ws.Range("M134").Value = -12790.845

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3000021.5
$ws.Range("I4").Value = 3000021.5
$ws.Range("K4").Value = 9000064.5
$ws.Range("M4").Value = -8999952.5
$ws.Range("H11").Value = 10153.429
$ws.Range("I11").Value = 10486.1
$ws.Range("J11").Value = 3500
$ws.Range("K11").Value = 31458.3
$ws.Range("L11").Value = 10500
$ws.Range("M11").Value = -31318.3
$ws.Range("N11").Value = -10780
$ws.Range("H60").Value = 1062.95
$ws.Range("I60").Value = 266.07693
$ws.Range("K60").Value = 798.2307900000001
$ws.Range("M60").Value = -547.2307900000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3376.25
$ws.Range("I80").Value = 3001.6667
$ws.Range("J80").Value = 4500
$ws.Range("K80").Value = 3001.6667
$ws.Range("L80").Value = 4500
$ws.Range("M80").Value = -2003.6667
$ws.Range("N80").Value = -6496
$ws.Range("H83").Value = 3376.25
$ws.Range("I83").Value = 3001.6667
$ws.Range("J83").Value = 4500
$ws.Range("K83").Value = 15008.3335
$ws.Range("L83").Value = 22500
$ws.Range("M83").Value = -10016.3335
$ws.Range("N83").Value = -32484
$ws.Range("H97").Value = 1365.625
$ws.Range("I97").Value = 1833.3334
$ws.Range("J97").Value = 1085
$ws.Range("K97").Value = 1833.3334
$ws.Range("L97").Value = 1085
$ws.Range("M97").Value = -1337.3334
$ws.Range("N97").Value = -2077
$ws.Range("H132").Value = 129181.625
$ws.Range("I132").Value = 204207.6
$ws.Range("J132").Value = 4138.3335
$ws.Range("K132").Value = 612622.8
$ws.Range("L132").Value = 12415.0005
$ws.Range("M132").Value = -610092.8
$ws.Range("N132").Value = -17475.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6220.75
$ws.Range("J46").Value = 6695.4546
$ws.Range("L46").Value = 6695.4546
$ws.Range("N46").Value = -7071.4546
$ws.Range("H55").Value = 1132.6666
$ws.Range("I55").Value = 898.9231
$ws.Range("J55").Value = 1512.5
$ws.Range("K55").Value = 898.9231
$ws.Range("L55").Value = 1512.5
$ws.Range("M55").Value = -725.9231
$ws.Range("N55").Value = -1858.5
$ws.Range("H74").Value = 50192
$ws.Range("I74").Value = 50192
$ws.Range("K74").Value = 50192
$ws.Range("M74").Value = -49194
$ws.Range("H77").Value = 50192
$ws.Range("I77").Value = 50192
$ws.Range("K77").Value = 150576
$ws.Range("M77").Value = -145584
$ws.Range("H132").Value = 17800.8
$ws.Range("I132").Value = 17800.8
$ws.Range("K132").Value = 53402.39999999999
$ws.Range("M132").Value = -50872.39999999999
$ws.Range("I136").Value = 2501.3333
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 7503.999899999999
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -4953.999899999999
$ws.Range("N136").Value = -12600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7916.5557
$ws.Range("I62").Value = 3083.3333
$ws.Range("K62").Value = 3083.3333
$ws.Range("M62").Value = -2459.3333
$ws.Range("H65").Value = 7916.5557
$ws.Range("I65").Value = 3083.3333
$ws.Range("K65").Value = 15416.6665
$ws.Range("M65").Value = -12296.6665
$ws.Range("H136").Value = 2951.8965
$ws.Range("J136").Value = 4554.4546
$ws.Range("L136").Value = 13663.3638
$ws.Range("N136").Value = -18763.3638

Write-Output "Applied all Halicarnassus_Profits updates"